$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.698.19"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "2.944.36"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'590.81"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'148.77"
$ws.Range("E6").Value = "  +7.46%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.508"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.942.53"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").Value = "'7.14"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +10.23%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "'0.0000238"
$ws.Range("E13").Value = "  +9.54%  "
$ws.Range("D14").Value = "'32.53"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "3.432.43"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "62.699.80"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "'6.65"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "2.943.37"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "'438.32"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'0.665"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  +7.86%  "
$ws.Range("D25").Value = "'80.38"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'11.91"
$ws.Range("E26").Value = "  +4.63%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "  +11.73%  "
$ws.Range("E30").Value = "  +23.98%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").Value = "'2.16"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("E33").Value = "  +5.93%  "
$ws.Range("D34").Value = "'26.18"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("E37").Value = "  +12.38%  "
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("D39").Value = "'49.62"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +7.42%  "
$ws.Range("D41").Value = "'8.41"
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D43").Value = "'0.279"
$ws.Range("E43").Value = "  +5.19%  "
$ws.Range("D44").Value = "'39.87"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("D46").Value = "2.705.48"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "'0.0341"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").Value = "'358.18"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "'22.79"
